$d = $word.ActiveDocument

# Locate the first paragraph ("This is a Microsoft word document.") and find
# the insertion point just after its text, i.e. before the paragraph mark.
$p1 = $d.Paragraphs(1)
$insertPos = $p1.Range.End - 1

# The three fragments to append, each of which must land in its own <w:r>
# run (matching the target XML: no run-property differences, just separate
# runs). Plain sequential InsertAfter calls get folded back into a single
# run because the neighbouring runs end up with identical formatting, so a
# temporary bookmark is dropped at each seam immediately before the text on
# that side of the seam is written; a bookmark edge blocks the run-merge
# there. The bookmarks are deleted again once all the text is in place, so
# they leave no trace in the saved document.
$fragments = @(" (", "Changed main", ")")

$pos = $insertPos
for ($i = 0; $i -lt $fragments.Length; $i++) {
    $seam = $d.Range($pos, $pos)
    $bookmarkName = "ironSeam" + $i
    $d.Bookmarks.Add($bookmarkName, $seam)

    $target = $d.Range($pos, $pos)
    $target.InsertAfter($fragments[$i])

    $pos = $pos + $fragments[$i].Length
}

for ($i = 0; $i -lt $fragments.Length; $i++) {
    $bookmarkName = "ironSeam" + $i
    $d.Bookmarks($bookmarkName).Delete()
}
